$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "319.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.59%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "48.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "15.38%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.251"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.69%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08085"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.31%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.594"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.10%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.642"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.20%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "25.24%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1293"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8.06%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1930"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.05%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09404"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.80%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04596"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "9.69%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.23%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001331"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "5.60%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04169"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.10%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005875"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.10%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.341"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.34%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.428"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.84%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3406"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.95%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.096"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.94%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.70%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.08%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001308"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.96%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004246"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "7.99%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001354"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.82%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003548"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.72%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02699"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "11.33%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05663"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.58%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006315"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.49%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007921"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.41%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1443"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.94%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007732"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.98%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "11.42%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.98%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006910"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "9.44%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.77%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05731"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-32.84%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004010"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.54%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.77%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.77%"
